$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.4767139967802052
$ws.Range("D2").Value = 0.1414110798342989
$ws.Range("E2").Value = 0.08888232657264439
$ws.Range("F2").Value = 2.690921233113954
$ws.Range("G2").Value = 0.002552488474298453
$ws.Range("I2").Value = 2.049678246306939
$ws.Range("K2").Value = 0.7548592873351652
$ws.Range("L2").Value = 0.2771730801130587
$ws.Range("M2").Value = 0.178628610453579
$ws.Range("B3").Value = 0.4706384337675473
$ws.Range("D3").Value = 0.1405149052111483
$ws.Range("E3").Value = 0.08616212275616597
$ws.Range("F3").Value = 2.624404347715441
$ws.Range("G3").Value = 0.002557501898061155
$ws.Range("I3").Value = 2.026523120186823
$ws.Range("K3").Value = 0.6910621540522186
$ws.Range("L3").Value = 0.2691312886122574
$ws.Range("M3").Value = 0.175018353982928
$ws.Range("B4").Value = 0.4672450760307214
$ws.Range("D4").Value = 0.1399665884082211
$ws.Range("E4").Value = 0.08446149064471875
$ws.Range("F4").Value = 2.584557025517398
$ws.Range("G4").Value = 0.002560742046508135
$ws.Range("I4").Value = 2.012720892699498
$ws.Range("K4").Value = 0.6524240016238139
$ws.Range("L4").Value = 0.2643577862079525
$ws.Range("M4").Value = 0.1729184682892466
$ws.Range("B5").Value = 0.465947247614622
$ws.Range("D5").Value = 0.1397436004083303
$ws.Range("E5").Value = 0.08376070041519945
$ws.Range("F5").Value = 2.568567719350256
$ws.Range("G5").Value = 0.002562103278566319
$ws.Range("I5").Value = 2.007200331690498
$ws.Range("K5").Value = 0.6368116784981055
$ws.Range("L5").Value = 0.2624537518054098
$ws.Range("M5").Value = 0.1720921504452058
$ws.Range("B6").Value = 0.4657368845127792
$ws.Range("D6").Value = 0.1397066000023202
$ws.Range("E6").Value = 0.08364386200724816
$ws.Range("F6").Value = 2.56592769859023
$ws.Range("G6").Value = 0.002562331780748653
$ws.Range("I6").Value = 2.006289913207134
$ws.Range("K6").Value = 0.6342272667437499
$ws.Range("L6").Value = 0.2621400748251688
$ws.Range("M6").Value = 0.1719567176306427
$ws.Range("B7").Value = 0.467227228617844
$ws.Range("D7").Value = 0.1399635793115124
$ws.Range("E7").Value = 0.08445207115482489
$ws.Range("F7").Value = 2.584340382395581
$ws.Range("G7").Value = 0.00256076023916121
$ws.Range("I7").Value = 2.012646020132948
$ws.Range("K7").Value = 0.6522129111858135
$ws.Range("L7").Value = 0.2643319409953477
$ws.Range("M7").Value = 0.1729072051929812
$ws.Range("B8").Value = 0.4745493218450463
$ws.Range("D8").Value = 0.1411016572628725
$ws.Range("E8").Value = 0.08795063982634055
$ws.Range("F8").Value = 2.66777907603003
$ws.Range("G8").Value = 0.002554183586386575
$ws.Range("I8").Value = 2.041607937512126
$ws.Range("K8").Value = 0.7327505947080226
$ws.Range("L8").Value = 0.2743661405856699
$ws.Range("M8").Value = 0.1773595600172158
$ws.Range("B9").Value = 0.4915737201678212
$ws.Range("D9").Value = 0.1433503223696491
$ws.Range("E9").Value = 0.09457565138248469
$ws.Range("F9").Value = 2.839357319464824
$ws.Range("G9").Value = 0.002542565004074581
$ws.Range("I9").Value = 2.101718734317998
$ws.Range("K9").Value = 0.8949784315719285
$ws.Range("L9").Value = 0.2953510423216414
$ws.Range("M9").Value = 0.1870172528434786
$ws.Range("B10").Value = 0.5056975715906447
$ws.Range("D10").Value = 0.1450149468403268
$ws.Range("E10").Value = 0.09930755363051347
$ws.Range("F10").Value = 2.970371687043553
$ws.Range("G10").Value = 0.00253479926099284
$ws.Range("I10").Value = 2.147940883511424
$ws.Range("K10").Value = 1.016881776255019
$ws.Range("L10").Value = 0.3115751311753741
$ws.Range("M10").Value = 0.1946784154479957
$ws.Range("B11").Value = 0.5124720669965939
$ws.Range("D11").Value = 0.1457754330656371
$ws.Range("E11").Value = 0.1014325944412349
$ws.Range("F11").Value = 3.031073563916237
$ws.Range("G11").Value = 0.002531431835749677
$ws.Range("I11").Value = 2.169423849867314
$ws.Range("K11").Value = 1.07295021509708
$ws.Range("L11").Value = 0.3191331511058877
$ws.Range("M11").Value = 0.1982867641133339
$ws.Range("B12").Value = 0.515087462515055
$ws.Range("D12").Value = 0.1460639112350179
$ws.Range("E12").Value = 0.1022334745310332
$ws.Range("F12").Value = 3.054220005437116
$ws.Range("G12").Value = 0.00253018029868856
$ws.Range("I12").Value = 2.177625060620869
$ws.Range("K12").Value = 1.094271701974776
$ws.Range("L12").Value = 0.3220208499579087
$ws.Range("M12").Value = 0.1996708748351068
$ws.Range("B13").Value = 0.5145219680755986
$ws.Range("D13").Value = 0.1460017596088647
$ws.Range("E13").Value = 0.1020611587623215
$ws.Range("F13").Value = 3.049227868802518
$ws.Range("G13").Value = 0.002530448790537325
$ws.Range("I13").Value = 2.175855836889227
$ws.Range("K13").Value = 1.089675733073307
$ws.Range("L13").Value = 0.3213977903117211
$ws.Range("M13").Value = 0.1993719945152392
$ws.Range("B14").Value = 0.5126862354764796
$ws.Range("D14").Value = 0.145799156159967
$ws.Range("E14").Value = 0.1014985592639945
$ws.Range("F14").Value = 3.03297462183869
$ws.Range("G14").Value = 0.002531328398041332
$ws.Range("I14").Value = 2.17009724144522
$ws.Range("K14").Value = 1.074702543723333
$ws.Range("L14").Value = 0.319370209356336
$ws.Range("M14").Value = 0.1984002809361129
$ws.Range("B15").Value = 0.5115683067866144
$ws.Range("D15").Value = 0.1456751215901164
$ws.Range("E15").Value = 0.1011534563960588
$ws.Range("F15").Value = 3.023039911229574
$ws.Range("G15").Value = 0.002531870257095318
$ws.Range("I15").Value = 2.166578554564069
$ws.Range("K15").Value = 1.065542738441934
$ws.Range("L15").Value = 0.3181315997288579
$ws.Range("M15").Value = 0.1978073839035552
$ws.Range("B16").Value = 0.5052618565031537
$ws.Range("D16").Value = 0.1449653152664112
$ws.Range("E16").Value = 0.0991681345312756
$ws.Range("F16").Value = 2.966426985018245
$ws.Range("G16").Value = 0.002535022644373143
$ws.Range("I16").Value = 2.146546142357678
$ws.Range("K16").Value = 1.013230037178062
$ws.Range("L16").Value = 0.3110847792344202
$ws.Range("M16").Value = 0.1944450800768394
$ws.Range("B17").Value = 0.5014824112559211
$ws.Range("D17").Value = 0.1445307251890497
$ws.Range("E17").Value = 0.09794324697768175
$ws.Range("F17").Value = 2.931980214491347
$ws.Range("G17").Value = 0.002536998764466298
$ws.Range("I17").Value = 2.13437410064671
$ws.Range("K17").Value = 0.9812960619150886
$ws.Range("L17").Value = 0.3068073426051683
$ws.Range("M17").Value = 0.1924139707540071
$ws.Range("B18").Value = 0.4993414820452529
$ws.Range("D18").Value = 0.1442810635112366
$ws.Range("E18").Value = 0.097236125299343
$ws.Range("F18").Value = 2.912271059076602
$ws.Range("G18").Value = 0.002538150938427338
$ws.Range("I18").Value = 2.127415983385731
$ws.Range("K18").Value = 0.9629861973849074
$ws.Range("L18").Value = 0.3043637845363918
$ws.Range("M18").Value = 0.1912573318524871
$ws.Range("B19").Value = 0.4986222587225768
$ws.Range("D19").Value = 0.1441965834093679
$ws.Range("E19").Value = 0.09699625573406578
$ws.Range("F19").Value = 2.905615649217594
$ws.Range("G19").Value = 0.002538543721324117
$ws.Range("I19").Value = 2.125067441701916
$ws.Range("K19").Value = 0.9567966651440258
$ws.Range("L19").Value = 0.3035393047128139
$ws.Range("M19").Value = 0.1908677068491755
$ws.Range("B20").Value = 0.5018813353979681
$ws.Range("D20").Value = 0.1445769564469046
$ws.Range("E20").Value = 0.09807390635919688
$ws.Range("F20").Value = 2.935636387910762
$ws.Range("G20").Value = 0.002536786793677276
$ws.Range("I20").Value = 2.135665389951185
$ws.Range("K20").Value = 0.9846895056711276
$ws.Range("L20").Value = 0.3072609528494894
$ws.Range("M20").Value = 0.1926289852195353
$ws.Range("B21").Value = 0.513224078227438
$ws.Range("D21").Value = 0.1458586519152405
$ws.Range("E21").Value = 0.1016639111768107
$ws.Range("F21").Value = 3.03774424217093
$ws.Range("G21").Value = 0.002531069395595836
$ws.Range("I21").Value = 2.171786883469352
$ws.Range("K21").Value = 1.07909809062636
$ws.Range("L21").Value = 0.3199650623240871
$ws.Range("M21").Value = 0.198685216311091
$ws.Range("B22").Value = 0.5209287975023642
$ws.Range("D22").Value = 0.1466992349606997
$ws.Range("E22").Value = 0.1039879437044142
$ws.Range("F22").Value = 3.105410637002365
$ws.Range("G22").Value = 0.002527470443069434
$ws.Range("I22").Value = 2.195779736207001
$ws.Range("K22").Value = 1.141322351028805
$ws.Range("L22").Value = 0.3284174353177889
$ws.Range("M22").Value = 0.2027465195508356
$ws.Range("B23").Value = 0.516790033999456
$ws.Range("D23").Value = 0.1462503219436826
$ws.Range("E23").Value = 0.1027495552583453
$ws.Range("F23").Value = 3.069209978525322
$ws.Range("G23").Value = 0.002529378714963628
$ws.Range("I23").Value = 2.182938883561761
$ws.Range("K23").Value = 1.108063823000521
$ws.Range("L23").Value = 0.3238925308462228
$ws.Range("M23").Value = 0.2005694869275345
$ws.Range("B24").Value = 0.5017008823378717
$ws.Range("D24").Value = 0.1445560547057383
$ws.Range("E24").Value = 0.0980148443523845
$ws.Range("F24").Value = 2.933983136861826
$ws.Range("G24").Value = 0.002536882575453481
$ws.Range("I24").Value = 2.135081474262805
$ws.Range("K24").Value = 0.9831551759696993
$ws.Range("L24").Value = 0.3070558270831611
$ws.Range("M24").Value = 0.192531742688189
$ws.Range("B25").Value = 0.4866837332501035
$ws.Range("D25").Value = 0.1427399591543903
$ws.Range("E25").Value = 0.09280780022775659
$ws.Range("F25").Value = 2.792078278583261
$ws.Range("G25").Value = 0.002545572209643331
$ws.Range("I25").Value = 2.085098304383948
$ws.Range("K25").Value = 0.8506218876936771
$ws.Range("L25").Value = 0.2895331155155674
$ws.Range("M25").Value = 0.1843053155590155
